$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question = "Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?"

# Row 20
$ws.Range("A20").Value = $question
$ws.Range("B20").Value = "no"
$ws.Range("C20").Value = "neutral"
$ws.Range("D20").Value = "I totally understand! But just to let you know, we have a special offer running — a 20% discount on all products today only! Would you like me to show you some options?"
$ws.Range("E20").Value = "2025-10-31 16:23:32"

# Row 21
$ws.Range("A21").Value = $question
$ws.Range("B21").Value = "No response"
$ws.Range("C21").Value = "neutral"
$ws.Range("D21").Value = "User rejected offer after persuasion"
$ws.Range("E21").Value = "2025-10-31 16:23:50"
